# Apply "template filtering, highlighting of mandatory columns" edit.
#
# Summary of the change:
#  - A new "DEFAULT" column is inserted at C, shifting the old C..K columns
#    (MANDATORY, TYPE, PATTERN, LOWER, UPPER, HEADER) one place to the right
#    (new D..L).
#  - The new DEFAULT column (C) gets a boolean flag for every data row,
#    mirroring the same true/false pattern the MANDATORY column used to show
#    before this edit.
#  - The worksheet view is re-zoomed/re-scrolled and a new active selection
#    is stored.
#  - Basic page setup (paper size / orientation) is defined for the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert a new column at C -----------------------------------------
# This shifts the old C..K columns one place to the right (D..L), carries
# over their column widths/styles automatically, and extends the sheet
# dimension/row spans to the new L column.
$ws.Columns("C").Insert()

# --- 2. Populate the new "DEFAULT" column ---------------------------------
$ws.Range("C1").Value = "DEFAULT"
$ws.Columns("C").ColumnWidth = 8

$ws.Range("C2").Value  = $true
$ws.Range("C3").Value  = $true
$ws.Range("C4").Value  = $true
$ws.Range("C5").Value  = $true
$ws.Range("C6").Value  = $true
$ws.Range("C7").Value  = $true
$ws.Range("C8").Value  = $true
$ws.Range("C9").Value  = $true
$ws.Range("C10").Value = $true
$ws.Range("C11").Value = $false
$ws.Range("C12").Value = $false
$ws.Range("C13").Value = $false
$ws.Range("C14").Value = $true
$ws.Range("C15").Value = $true

# --- 3. Update the sheet view: zoom + active selection --------------------
$excel.ActiveWindow.Zoom = 127
$ws.Range("C16").Select()

# --- 4. Basic page setup ---------------------------------------------------
$ps = $ws.PageSetup
$ps.PaperSize = 9      # xlPaperA4
$ps.Orientation = 1    # xlPortrait
